# "Generate Report for Archive"
# The localization status report is regenerated: the in-flight file's
# status moves from "Ready for handoff" to "In Translation" on every
# sheet that surfaces it (Overview, zh-cn, de-de), and the Status /
# per-language columns that held the old, longer text are re-fit to
# the new, shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-language status lives in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

# --- zh-cn sheet: Status lives in column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null

# --- de-de sheet: Status lives in column C
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
